# Updated cryptos list on Sat Jun  8 02:21:32 UTC 2024 with GitHub Actions
# Refresh the scraped coin table (Coin/Link/Price/Volume columns) with the
# latest values. Cell values are prefixed with a leading apostrophe so
# Excel stores them as literal text (matching the source inline strings)
# instead of auto-converting number-looking text like "1.00" or "0.148"
# into numeric values; Style is reset to "Normal" right after so no new
# quote-prefix formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'69.407.87"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -1.90%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'3.688.73"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +0.06%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'682.19"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -3.04%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'162.70"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -4.13%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'3.688.55"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  -2.93%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -4.19%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'0.148"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -7.61%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'7.33"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -3.31%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  -2.49%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'0.0000239"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -4.75%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'33.60"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -6.16%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'4.313.37"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -2.91%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'3.690.58"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -3.44%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'69.434.52"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -1.90%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -1.44%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'16.28"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -5.72%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'6.64"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -6.42%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'483.19"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -2.18%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'9.91"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -7.22%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'0.667"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  -7.81%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'80.36"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -4.54%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'3.835.78"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -2.98%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'0.0000130"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -8.80%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 2).Value = "'InternetComputer(DFINITY)"
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "'11.54"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -4.22%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 2).Value = "'Dai"
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).Value = "'1.00"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +0.05%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'9.66"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -6.44%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'1.84"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -8.79%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -9.26%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'2.09"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -8.19%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'6.84"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  -6.50%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 2).Value = "'EthereumClassic"
$ws.Cells.Item(34, 2).Style = "Normal"
$ws.Cells.Item(34, 3).Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'27.11"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -6.54%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 2).Value = "'Kaspa"
$ws.Cells.Item(35, 2).Style = "Normal"
$ws.Cells.Item(35, 3).Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'0.167"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -5.34%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  +0.03%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "'3.657.76"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -3.10%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'8.53"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -5.71%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'6.36"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +7.26%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'0.0935"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -7.12%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 2).Value = "'Stacks"
$ws.Cells.Item(41, 2).Style = "Normal"
$ws.Cells.Item(41, 3).Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'2.26"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -3.94%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 2).Value = "'USDe"
$ws.Cells.Item(42, 2).Style = "Normal"
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'1.00"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  -0.01%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -0.08%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  -6.56%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'159.57"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -4.39%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 2).Value = "'OKB"
$ws.Cells.Item(46, 2).Style = "Normal"
$ws.Cells.Item(46, 3).Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(46, 3).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'48.51"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -0.91%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 2).Value = "'dogwifhat"
$ws.Cells.Item(47, 2).Style = "Normal"
$ws.Cells.Item(47, 3).Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(47, 3).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'2.90"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -10.48%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 2).Value = "'ONDO"
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'1.37"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +1.84%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 2).Value = "'FLOKI"
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'0.000290"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -9.35%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'29.91"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +5.71%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'395.70"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -7.19%  "
$ws.Cells.Item(51, 5).Style = "Normal"
